# "order agnostic binary search and leetcode easy arrays questions"
#
# Adds two new LeetCode practice rows (Rotate Image / Add to Array-Form of
# Integer) right after the existing three question rows, plus a classroom
# assignment note further down the sheet, mirroring the shape of the other
# question rows (Question Heading col A, hyperlinked Question Link col B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New question rows (5 & 6) - heading + link, same layout as rows 2-4.
$ws.Range("A5").Value = "Rotate Image"
$ws.Range("B5").Value = "https://leetcode.com/problems/rotate-image/"

$ws.Range("A6").Value = "Add to Array-Form of Integer"
$ws.Range("B6").Value = "https://leetcode.com/problems/add-to-array-form-of-integer/"

# Classroom assignment note further down the sheet.
$ws.Range("A19").Value = "Arrays Medium Questions Community Classroom Assignment"

# Turn the two new Question Link cells into real hyperlinks.
$ws.Hyperlinks.Add($ws.Range("B5"), "https://leetcode.com/problems/rotate-image/")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://leetcode.com/problems/add-to-array-form-of-integer/")

# Match the existing hyperlink-cell formatting (blue text) used by B2:B4,
# rather than the default Excel "Hyperlink" style the Add() call applies.
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)

# Leave the selection where the author left it when saving.
$ws.Range("A21").Select()
